$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like numbers need a temporary Text
# number format so Excel keeps them as strings (matching the source
# inline-string cells) instead of auto-converting to numeric values.
$numericLookingRefs = @("D5", "D6", "D11", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D29", "D30", "D31", "D32", "D37", "D38", "D41", "D42", "D43", "D47", "D51")
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "68.972.88"
$ws.Range("E2").Value2 = "  +0.89%  "
$ws.Range("D3").Value2 = "2.742.80"
$ws.Range("E3").Value2 = "  +3.52%  "
$ws.Range("E4").Value2 = "  +0.04%  "
$ws.Range("D5").Value2 = "605.24"
$ws.Range("E5").Value2 = "  +1.17%  "
$ws.Range("D6").Value2 = "166.78"
$ws.Range("E6").Value2 = "  +4.62%  "
$ws.Range("E8").Value2 = "  +0.46%  "
$ws.Range("D9").Value2 = "2.742.53"
$ws.Range("E9").Value2 = "  +3.56%  "
$ws.Range("E10").Value2 = "  -1.96%  "
$ws.Range("B11").Value2 = "Toncoin"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value2 = "5.36"
$ws.Range("E11").Value2 = "  +1.69%  "
$ws.Range("B12").Value2 = "Cardano"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value2 = "0.365"
$ws.Range("E12").Value2 = "  +3.78%  "
$ws.Range("E13").Value2 = "  -0.30%  "
$ws.Range("D14").Value2 = "28.91"
$ws.Range("E14").Value2 = "  +2.97%  "
$ws.Range("D15").Value2 = "3.247.30"
$ws.Range("E16").Value2 = "  -0.36%  "
$ws.Range("D17").Value2 = "68.928.12"
$ws.Range("E17").Value2 = "  +1.01%  "
$ws.Range("D18").Value2 = "2.734.26"
$ws.Range("E18").Value2 = "  +2.53%  "
$ws.Range("D19").Value2 = "11.95"
$ws.Range("E19").Value2 = "  +4.67%  "
$ws.Range("D20").Value2 = "7.75"
$ws.Range("E20").Value2 = "  +5.80%  "
$ws.Range("D21").Value2 = "369.13"
$ws.Range("E21").Value2 = "  +1.51%  "
$ws.Range("D22").Value2 = "4.57"
$ws.Range("E22").Value2 = "  +3.04%  "
$ws.Range("D23").Value2 = "4.97"
$ws.Range("E23").Value2 = "  +2.91%  "
$ws.Range("E24").Value2 = "  +3.13%  "
$ws.Range("D25").Value2 = "74.23"
$ws.Range("E25").Value2 = "  -0.99%  "
$ws.Range("E26").Value2 = "  -0.01%  "
$ws.Range("D27").Value2 = "9.98"
$ws.Range("E27").Value2 = "  +2.73%  "
$ws.Range("D29").Value2 = "0.0000106"
$ws.Range("E29").Value2 = "  +1.33%  "
$ws.Range("D30").Value2 = "600.12"
$ws.Range("E30").Value2 = "  +7.50%  "
$ws.Range("D31").Value2 = "1.00"
$ws.Range("E31").Value2 = "  -3.72%  "
$ws.Range("D32").Value2 = "8.34"
$ws.Range("E32").Value2 = "  +3.82%  "
$ws.Range("E33").Value2 = "  +3.68%  "
$ws.Range("E34").Value2 = "  +5.96%  "
$ws.Range("E35").Value2 = "  +3.20%  "
$ws.Range("E36").Value2 = "  +4.02%  "
$ws.Range("D37").Value2 = "1.00"
$ws.Range("E37").Value2 = "  +0.03%  "
$ws.Range("D38").Value2 = "163.42"
$ws.Range("E38").Value2 = "  +2.31%  "
$ws.Range("E39").Value2 = "  +1.41%  "
$ws.Range("E40").Value2 = "  +3.49%  "
$ws.Range("D41").Value2 = "1.92"
$ws.Range("E41").Value2 = "  +2.27%  "
$ws.Range("D42").Value2 = "5.51"
$ws.Range("E42").Value2 = "  +2.64%  "
$ws.Range("D43").Value2 = "2.71"
$ws.Range("E43").Value2 = "  +2.59%  "
$ws.Range("D45").Value2 = "0.0₆0318"
$ws.Range("E45").Value2 = "  -5.12%  "
$ws.Range("E46").Value2 = "  +0.01%  "
$ws.Range("D47").Value2 = "158.73"
$ws.Range("E47").Value2 = "  +0.45%  "
$ws.Range("E48").Value2 = "  +5.36%  "
$ws.Range("E49").Value2 = "  +6.65%  "
$ws.Range("E50").Value2 = "  +7.74%  "
$ws.Range("D51").Value2 = "22.12"
$ws.Range("E51").Value2 = "  -0.75%  "

# Restore the original (default) cell style now that the text value is set,
# so we don't leave a stray number-format style on these cells.
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).Style = "Normal"
}
